# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F column) counts to the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (1st sheet) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value  = 92
$wsExpo.Range("F7").Value  = 1315
$wsExpo.Range("F13").Value = 190
$wsExpo.Range("F20").Value = 294
$wsExpo.Range("F21").Value = 3240
$wsExpo.Range("F24").Value = 941
$wsExpo.Range("F28").Value = 1667
$wsExpo.Range("F35").Value = 913
$wsExpo.Range("F36").Value = 1949
$wsExpo.Range("F38").Value = 366
$wsExpo.Range("F39").Value = 112
$wsExpo.Range("F42").Value = 916
$wsExpo.Range("F43").Value = 819
$wsExpo.Range("F47").Value = 294
$wsExpo.Range("F49").Value = 3385

# --- Sheet "全部类型" (4th sheet) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 92
$wsAll.Range("F8").Value  = 1315
$wsAll.Range("F14").Value = 190
$wsAll.Range("F20").Value = 294
$wsAll.Range("F21").Value = 3240
$wsAll.Range("F27").Value = 1667
$wsAll.Range("F33").Value = 1940
$wsAll.Range("F35").Value = 913
$wsAll.Range("F36").Value = 1949
$wsAll.Range("F37").Value = 366
$wsAll.Range("F38").Value = 112
$wsAll.Range("F40").Value = 916
$wsAll.Range("F41").Value = 819
$wsAll.Range("F45").Value = 294
$wsAll.Range("F48").Value = 3385
